# Auto-generated edit script: updates profit-calculation columns (H:N)
# on specific rows across all 8 item-category worksheets, per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 490.58823
$ws.Range("I19").Value = 440
$ws.Range("J19").Value = 511.66666
$ws.Range("K19").Value = 440
$ws.Range("L19").Value = 511.66666
$ws.Range("M19").Value = -265
$ws.Range("N19").Value = -861.66666

$ws.Range("H138").Value = 2364.456
$ws.Range("I138").Value = 2710.3333
$ws.Range("J138").Value = 2272.2222
$ws.Range("K138").Value = 8130.999899999999
$ws.Range("L138").Value = 6816.6666
$ws.Range("M138").Value = -2990.999899999999
$ws.Range("N138").Value = -17096.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 26432.666
$ws.Range("I44").Value = 10000
$ws.Range("J44").Value = 34649
$ws.Range("K44").Value = 10000
$ws.Range("L44").Value = 34649
$ws.Range("M44").Value = -9512
$ws.Range("N44").Value = -35625

$ws.Range("H45").Value = 4306.1177
$ws.Range("I45").Value = 4002.2222
$ws.Range("J45").Value = 4648
$ws.Range("K45").Value = 4002.2222
$ws.Range("L45").Value = 4648
$ws.Range("M45").Value = -3625.2222
$ws.Range("N45").Value = -5402

$ws.Range("H55").Value = 23950.666
$ws.Range("J55").Value = 23950.666
$ws.Range("L55").Value = 23950.666
$ws.Range("N55").Value = -24580.666

$ws.Range("H61").Value = 1602.3715
$ws.Range("I61").Value = 1803.6086
$ws.Range("J61").Value = 1216.6666
$ws.Range("K61").Value = 1803.6086
$ws.Range("L61").Value = 1216.6666
$ws.Range("M61").Value = -1591.6086
$ws.Range("N61").Value = -1640.6666

$ws.Range("H80").Value = 26250.8
$ws.Range("J80").Value = 26250.8
$ws.Range("L80").Value = 26250.8
$ws.Range("N80").Value = -28246.8

$ws.Range("H83").Value = 26250.8
$ws.Range("J83").Value = 26250.8
$ws.Range("L83").Value = 78752.39999999999
$ws.Range("N83").Value = -88736.39999999999

$ws.Range("H102").Value = 125000960
$ws.Range("I102").Value = 1056.6666
$ws.Range("J102").Value = 500000670
$ws.Range("K102").Value = 1056.6666
$ws.Range("L102").Value = 500000670
$ws.Range("M102").Value = 565.3334
$ws.Range("N102").Value = -500003914

$ws.Range("H122").Value = 1548.7142
$ws.Range("I122").Value = 1723.3
$ws.Range("J122").Value = 1390
$ws.Range("K122").Value = 5169.9
$ws.Range("L122").Value = 4170
$ws.Range("M122").Value = -2719.9
$ws.Range("N122").Value = -9070

$ws.Range("H135").Value = 36765.5
$ws.Range("J135").Value = 36765.5
$ws.Range("L135").Value = 36765.5
$ws.Range("N135").Value = -46905.5

$ws.Range("H136").Value = 1602.3715
$ws.Range("I136").Value = 1803.6086
$ws.Range("J136").Value = 1216.6666
$ws.Range("K136").Value = 5410.825800000001
$ws.Range("L136").Value = 3649.9998
$ws.Range("M136").Value = -2860.825800000001
$ws.Range("N136").Value = -8749.9998

$ws.Range("H139").Value = 100488.89
$ws.Range("J139").Value = 100488.89
$ws.Range("L139").Value = 100488.89
$ws.Range("N139").Value = -110768.89

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 34887
$ws.Range("J35").Value = 34887
$ws.Range("L35").Value = 34887
$ws.Range("N35").Value = -35507

$ws.Range("H81").Value = 30000
$ws.Range("J81").Value = 30000
$ws.Range("L81").Value = 30000
$ws.Range("N81").Value = -32122

$ws.Range("H82").Value = 28511.941
$ws.Range("J82").Value = 32151.857
$ws.Range("L82").Value = 32151.857
$ws.Range("N82").Value = -32917.857

$ws.Range("H84").Value = 30000
$ws.Range("J84").Value = 30000
$ws.Range("L84").Value = 90000
$ws.Range("N84").Value = -100608

$ws.Range("H85").Value = 28511.941
$ws.Range("J85").Value = 32151.857
$ws.Range("L85").Value = 32151.857
$ws.Range("N85").Value = -34803.857

$ws.Range("H107").Value = 3981.575
$ws.Range("I107").Value = 3632.7354
$ws.Range("J107").Value = 5958.3335
$ws.Range("K107").Value = 3632.7354
$ws.Range("L107").Value = 5958.3335
$ws.Range("M107").Value = -1712.7354
$ws.Range("N107").Value = -9798.333500000001

$ws.Range("H135").Value = 58834.707
$ws.Range("J135").Value = 58834.707
$ws.Range("L135").Value = 58834.707
$ws.Range("N135").Value = -68974.70699999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2709.1897
$ws.Range("I31").Value = 1857.7354
$ws.Range("J31").Value = 3915.4167
$ws.Range("K31").Value = 1857.7354
$ws.Range("L31").Value = 3915.4167
$ws.Range("M31").Value = -1562.7354
$ws.Range("N31").Value = -4505.4167

$ws.Range("H34").Value = 2709.1897
$ws.Range("I34").Value = 1857.7354
$ws.Range("J34").Value = 3915.4167
$ws.Range("K34").Value = 1857.7354
$ws.Range("L34").Value = 3915.4167
$ws.Range("M34").Value = -1655.7354
$ws.Range("N34").Value = -4319.4167

$ws.Range("H51").Value = 9258.4
$ws.Range("J51").Value = 9258.4
$ws.Range("L51").Value = 9258.4
$ws.Range("N51").Value = -10730.4

$ws.Range("H58").Value = 1685.1578
$ws.Range("I58").Value = 1775.25
$ws.Range("J58").Value = 1204.6666
$ws.Range("K58").Value = 1775.25
$ws.Range("L58").Value = 1204.6666
$ws.Range("M58").Value = -1572.25
$ws.Range("N58").Value = -1610.6666

$ws.Range("H60").Value = 36938.285
$ws.Range("J60").Value = 36938.285
$ws.Range("L60").Value = 36938.285
$ws.Range("N60").Value = -37960.285

$ws.Range("H61").Value = 9258.4
$ws.Range("J61").Value = 9258.4
$ws.Range("L61").Value = 9258.4
$ws.Range("N61").Value = -9954.4

$ws.Range("H68").Value = 17326.666
$ws.Range("J68").Value = 17326.666
$ws.Range("L68").Value = 17326.666
$ws.Range("N68").Value = -18824.666

$ws.Range("H71").Value = 17326.666
$ws.Range("J71").Value = 17326.666
$ws.Range("L71").Value = 51979.99800000001
$ws.Range("N71").Value = -59467.99800000001

$ws.Range("H109").Value = 10885.714
$ws.Range("J109").Value = 10885.714
$ws.Range("L109").Value = 10885.714
$ws.Range("N109").Value = -12965.714

$ws.Range("H132").Value = 2783.3333
$ws.Range("I132").Value = 2756.25
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 8268.75
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -5738.75
$ws.Range("N132").Value = -14060

$ws.Range("H136").Value = 1685.1578
$ws.Range("I136").Value = 1775.25
$ws.Range("J136").Value = 1204.6666
$ws.Range("K136").Value = 5325.75
$ws.Range("L136").Value = 3613.9998
$ws.Range("M136").Value = -2775.75
$ws.Range("N136").Value = -8713.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 551.0833
$ws.Range("I44").Value = 192.5
$ws.Range("J44").Value = 563.4483
$ws.Range("K44").Value = 577.5
$ws.Range("L44").Value = 1690.3449
$ws.Range("M44").Value = -179.5
$ws.Range("N44").Value = -2486.3449

$ws.Range("H103").Value = 1243.75
$ws.Range("J103").Value = 1500
$ws.Range("L103").Value = 4500
$ws.Range("N103").Value = -6258

$ws.Range("H123").Value = 3340
$ws.Range("I123").Value = 1865.7142
$ws.Range("K123").Value = 5597.142599999999
$ws.Range("M123").Value = -3147.142599999999

$ws.Range("H131").Value = 906.3570999999999
$ws.Range("I131").Value = 594.875
$ws.Range("J131").Value = 934.04443
$ws.Range("K131").Value = 1784.625
$ws.Range("L131").Value = 2802.13329
$ws.Range("M131").Value = 3255.375
$ws.Range("N131").Value = -12882.13329

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 15061.2
$ws.Range("J51").Value = 15061.2
$ws.Range("L51").Value = 15061.2
$ws.Range("N51").Value = -16079.2

$ws.Range("H102").Value = 2393.5
$ws.Range("I102").Value = 1191.3334
$ws.Range("J102").Value = 6000
$ws.Range("K102").Value = 1191.3334
$ws.Range("L102").Value = 6000
$ws.Range("M102").Value = 430.6666
$ws.Range("N102").Value = -9244

$ws.Range("H122").Value = 1801.0714
$ws.Range("I122").Value = 2084.5
$ws.Range("J122").Value = 1588.5
$ws.Range("K122").Value = 6253.5
$ws.Range("L122").Value = 4765.5
$ws.Range("M122").Value = -3803.5
$ws.Range("N122").Value = -9665.5

$ws.Range("H123").Value = 29170.4
$ws.Range("J123").Value = 29170.4
$ws.Range("L123").Value = 29170.4
$ws.Range("N123").Value = -34070.4

$ws.Range("H126").Value = 12092.25
$ws.Range("I126").Value = 3231
$ws.Range("J126").Value = 19772
$ws.Range("K126").Value = 9693
$ws.Range("L126").Value = 59316
$ws.Range("M126").Value = -7223
$ws.Range("N126").Value = -64256

$ws.Range("H132").Value = 2420.8857
$ws.Range("J132").Value = 2796.7222
$ws.Range("L132").Value = 8390.1666
$ws.Range("N132").Value = -13450.1666

$ws.Range("H140").Value = 42693.332
$ws.Range("J140").Value = 42693.332
$ws.Range("L140").Value = 42693.332
$ws.Range("N140").Value = -53053.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2885.7144
$ws.Range("I61").Value = 1825
$ws.Range("K61").Value = 1825
$ws.Range("M61").Value = -1623

$ws.Range("H109").Value = 22142.5
$ws.Range("J109").Value = 22142.5
$ws.Range("L109").Value = 22142.5
$ws.Range("N109").Value = -24916.5

$ws.Range("H113").Value = 2885.7144
$ws.Range("I113").Value = 1825
$ws.Range("K113").Value = 1825
$ws.Range("M113").Value = 345

$ws.Range("H122").Value = 5143.4287
$ws.Range("I122").Value = 2502
$ws.Range("J122").Value = 6200
$ws.Range("K122").Value = 7506
$ws.Range("L122").Value = 18600
$ws.Range("M122").Value = -5056
$ws.Range("N122").Value = -23500

$ws.Range("H132").Value = 3378.325
$ws.Range("I132").Value = 3159.5
$ws.Range("J132").Value = 4034.8
$ws.Range("K132").Value = 9478.5
$ws.Range("L132").Value = 12104.4
$ws.Range("M132").Value = -6948.5
$ws.Range("N132").Value = -17164.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 17389
$ws.Range("J109").Value = 17389
$ws.Range("L109").Value = 17389
$ws.Range("N109").Value = -20163

$ws.Range("H122").Value = 7866.657
$ws.Range("I122").Value = 11435.81
$ws.Range("J122").Value = 2512.9285
$ws.Range("K122").Value = 34307.43
$ws.Range("L122").Value = 7538.7855
$ws.Range("M122").Value = -31857.43
$ws.Range("N122").Value = -12438.7855
